$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to stay a plain text string (avoid Excel's automatic
    # date-like conversion of values such as "01/04/2025").
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# New region order (rows 2-7) and refreshed quarter / value figures (rows 2-10)
$ws.Range("A2").Value = "Santa Catarina"
Set-TextValue $ws.Range("C2") "01/04/2025"
$ws.Range("D2").Value = 97.75

$ws.Range("A3").Value = "Rondônia"
Set-TextValue $ws.Range("C3") "01/04/2025"
$ws.Range("D3").Value = 97.69

$ws.Range("A4").Value = "Mato Grosso"
Set-TextValue $ws.Range("C4") "01/04/2025"
$ws.Range("D4").Value = 97.21

$ws.Range("A5").Value = "Mato Grosso do Sul"
Set-TextValue $ws.Range("C5") "01/04/2025"
$ws.Range("D5").Value = 97.15

$ws.Range("A6").Value = "Espírito Santo"
Set-TextValue $ws.Range("C6") "01/04/2025"
$ws.Range("D6").Value = 96.86

$ws.Range("A7").Value = "Paraná"
Set-TextValue $ws.Range("C7") "01/04/2025"
$ws.Range("D7").Value = 96.19

Set-TextValue $ws.Range("C8") "01/04/2025"
$ws.Range("D8").Value = 91.95

Set-TextValue $ws.Range("C9") "01/04/2025"
$ws.Range("D9").Value = 94.24

Set-TextValue $ws.Range("C10") "01/04/2025"
$ws.Range("D10").Value = 91.77
